# Fruta / hortaliza, semanal
# This workbook's data rows (2-24) get their mutable fields
# (Fecha, Calidad, Volumen, Precio minimo, Precio maximo,
#  Precio promedio ponderado, Origen, Precio $/Kg) reshuffled
# between rows. Columns A,B,C,E,F,G,H,N,Q,R are identical
# across all rows and stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row, as derived from the diff (row -> values)
$rows = @{
    2  = @{ D = 44211; I = "Primera"; J = 28;  K = 8000; L = 8500; M = 8214; O = "Región Metropolitana";   P = 821 }
    3  = @{ D = 44812; I = "Primera"; J = 60;  K = 7000; L = 8000; M = 7500; O = "Provincia de Diguillín"; P = 750 }
    4  = @{ D = 44798; I = "Primera"; J = 80;  K = 7000; L = 7000; M = 7000; O = "Provincia de Diguillín"; P = 700 }
    5  = @{ D = 44838; I = "Primera"; J = 120; K = 6500; L = 7000; M = 6750; O = "Provincia de Diguillín"; P = 675 }
    6  = @{ D = 44791; I = "Primera"; J = 100; K = 8500; L = 9000; M = 8750; O = "Región Metropolitana";   P = 875 }
    7  = @{ D = 44775; I = "Primera"; J = 60;  K = 8000; L = 8000; M = 8000; O = "Región Metropolitana";   P = 800 }
    8  = @{ D = 44799; I = "Primera"; J = 60;  K = 7000; L = 7000; M = 7000; O = "Provincia de Diguillín"; P = 700 }
    9  = @{ D = 44831; I = "Primera"; J = 60;  K = 7000; L = 7500; M = 7250; O = "Provincia de Diguillín"; P = 725 }
    10 = @{ D = 44784; I = "Primera"; J = 100; K = 8000; L = 9000; M = 8500; O = "Región Metropolitana";   P = 850 }
    11 = @{ D = 44806; I = "Primera"; J = 120; K = 7000; L = 7500; M = 7250; O = "Provincia de Diguillín"; P = 725 }
    12 = @{ D = 44790; I = "Primera"; J = 60;  K = 8500; L = 9000; M = 8750; O = "Región Metropolitana";   P = 875 }
    13 = @{ D = 44813; I = "Primera"; J = 120; K = 7000; L = 7500; M = 7250; O = "Provincia de Diguillín"; P = 725 }
    14 = @{ D = 44846; I = "Primera"; J = 100; K = 6500; L = 7000; M = 6750; O = "Provincia de Diguillín"; P = 675 }
    15 = @{ D = 44841; I = "Primera"; J = 60;  K = 6500; L = 7000; M = 6750; O = "Provincia de Diguillín"; P = 675 }
    16 = @{ D = 44817; I = "Primera"; J = 60;  K = 7000; L = 7000; M = 7000; O = "Provincia de Diguillín"; P = 700 }
    17 = @{ D = 44817; I = "Segunda"; J = 60;  K = 8000; L = 8000; M = 8000; O = "Provincia de Diguillín"; P = 800 }
    18 = @{ D = 44819; I = "Primera"; J = 100; K = 7000; L = 8000; M = 7500; O = "Provincia de Diguillín"; P = 750 }
    19 = @{ D = 44810; I = "Primera"; J = 60;  K = 7000; L = 8000; M = 7500; O = "Provincia de Diguillín"; P = 750 }
    20 = @{ D = 44782; I = "Primera"; J = 120; K = 8000; L = 9000; M = 8500; O = "Región Metropolitana";   P = 850 }
    21 = @{ D = 44847; I = "Primera"; J = 100; K = 6500; L = 7000; M = 6750; O = "Provincia de Diguillín"; P = 675 }
    22 = @{ D = 44804; I = "Primera"; J = 80;  K = 7000; L = 7500; M = 7250; O = "Provincia de Diguillín"; P = 725 }
    23 = @{ D = 44203; I = "Primera"; J = 27;  K = 7000; L = 8000; M = 7556; O = "Región Metropolitana";   P = 756 }
    24 = @{ D = 44980; I = "Primera"; J = 60;  K = 7500; L = 8000; M = 7750; O = "Provincia de Diguillín"; P = 775 }
}

foreach ($r in $rows.Keys) {
    $v = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $v.D   # D: Fecha
    $ws.Cells.Item($r, 9).Value  = $v.I   # I: Calidad
    $ws.Cells.Item($r, 10).Value = $v.J   # J: Volumen
    $ws.Cells.Item($r, 11).Value = $v.K   # K: Precio minimo
    $ws.Cells.Item($r, 12).Value = $v.L   # L: Precio maximo
    $ws.Cells.Item($r, 13).Value = $v.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 15).Value = $v.O   # O: Origen
    $ws.Cells.Item($r, 16).Value = $v.P   # P: Precio $/Kg
}
